# Insert a new price record as row 280, pushing the existing rows
# 280:296 down to 281:297 (used range grows from A1:T296 to A1:T297).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("280:280").Insert()

$ws.Range("A280").Value = 4
$ws.Range("B280").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C280").Value = "Los Lagos"
$ws.Range("D280").Value = 44931
$ws.Range("E280").Value = 10
$ws.Range("F280").Value = "Fruta"
$ws.Range("G280").Value = 100108
$ws.Range("H280").Value = "Tropicales y subtropicales"
$ws.Range("I280").Value = 100108002
$ws.Range("J280").Value = "Mango"
$ws.Range("K280").Value = "Sin especificar"
$ws.Range("L280").Value = "Primera"
$ws.Range("M280").Value = 120
$ws.Range("N280").Value = 11000
$ws.Range("O280").Value = 11000
$ws.Range("P280").Value = 11000
$ws.Range("Q280").Value = "$/bandeja 4 kilos"
$ws.Range("R280").Value = "Brasil"
$ws.Range("S280").Value = 2750
$ws.Range("T280").Value = 4
